# Restricted and view users cannot access Delete/DeleteConfirmed routes.
# Append four new Test Plan rows documenting that restricted/view users
# are redirected away from the Quizzes Delete and DeleteConfirmed routes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing values/strings already present in the sheet so the
# shared-string table points at the same entries as the similar
# "Quizzes/Create" navigation rows above (column B = "Navigation",
# column D = "User is redirected to /Error/AccessDenied",
# column E = "Same as expected", column F = "Pass").
$navigation = $ws.Range("B19").Value2
$redirected = $ws.Range("D19").Value2
$sameAsExpected = $ws.Range("E19").Value2
$pass = $ws.Range("F19").Value2

# Row 30
$ws.Range("A30").Value = 27
$ws.Range("B30").Value = $navigation
$ws.Range("C30").Value = "Restricted user navigates to /Quizzes/Delete"
$ws.Range("D30").Value = $redirected
$ws.Range("E30").Value = $sameAsExpected
$ws.Range("F30").Value = $pass

# Row 31
$ws.Range("A31").Value = 28
$ws.Range("B31").Value = $navigation
$ws.Range("C31").Value = "View user navigates to /Quizzes/Delete"
$ws.Range("D31").Value = $redirected
$ws.Range("E31").Value = $sameAsExpected
$ws.Range("F31").Value = $pass

# Row 32
$ws.Range("A32").Value = 29
$ws.Range("B32").Value = $navigation
$ws.Range("C32").Value = "Restricted user navigates to /Quizzes/DeleteConfirmed"
$ws.Range("D32").Value = $redirected
$ws.Range("E32").Value = $sameAsExpected
$ws.Range("F32").Value = $pass

# Row 33
$ws.Range("A33").Value = 30
$ws.Range("B33").Value = $navigation
$ws.Range("C33").Value = "View user navigates to /Quizzes/DeleteConfirmed"
$ws.Range("D33").Value = $redirected
$ws.Range("E33").Value = $sameAsExpected
$ws.Range("F33").Value = $pass

# Match the formatting of the other "Scenario"/"Expected Results" columns
# (wrapped text, same style as the rest of column C/D).
$ws.Range("C30:D33").WrapText = $true

# Wrapped long scenario text spans two lines for rows 30, 32 and 33 (row 31's
# scenario text is short enough to fit on one line), matching row heights
# used elsewhere in the sheet for two-line wrapped cells.
$ws.Rows(30).RowHeight = 29.15
$ws.Rows(32).RowHeight = 29.15
$ws.Rows(33).RowHeight = 29.15

# Move the selection down to the newly-added last row/cell, consistent with
# having just finished entering this data at the bottom of the sheet.
$ws.Range("F33").Select()
